$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 100, shifting existing rows 100:122 down to 101:123
$ws.Rows.Item(100).Insert()

# Copy static columns A:K from row 101 (the row that was pushed down, previously row 100)
# into the new row 100, since they share identical values across this data block.
$ws.Range("A101:K101").Copy()
$ws.Range("A100").PasteSpecial()

# Also copy the date cell style (numeric date format) from row 101's D cell to D100
$ws.Range("D101").Copy()
$ws.Range("D100").PasteSpecial(-4122)

# Now set the new row's specific values
$ws.Range("D100").Value = 44995
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 400
$ws.Range("N100").Value = 17000
$ws.Range("O100").Value = 18000
$ws.Range("P100").Value = 17500
$ws.Range("Q100").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R100").Value = "Región de O'Higgins"
$ws.Range("S100").Value = 972
$ws.Range("T100").Value = 18
